# Update "想去人数" (want-to-go count) figures for four events that appear
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1048
$wsExhibit.Range("F3").Value = 302
$wsExhibit.Range("F4").Value = 2794
$wsExhibit.Range("F6").Value = 596

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1048
$wsAll.Range("F5").Value = 302
$wsAll.Range("F6").Value = 2794
$wsAll.Range("F8").Value = 596
